# Scheduled-runner price/profit refresh for the Leve profit tracker.
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N) for a
# set of leves across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets
# to reflect newly pulled market-board data. Some rows previously had no
# HQ profit (column N) because HQ price data was unavailable; those now
# gain a value, while a few rows that previously had stale HQ profit
# figures have lost their market data and had that cell cleared.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1343
$ws.Range("I86").Value = 827
$ws.Range("J86").Value = 1755.8
$ws.Range("K86").Value = 827
$ws.Range("L86").Value = 1755.8
$ws.Range("M86").Value = 296
$ws.Range("N86").Value = -4001.8
$ws.Range("H89").Value = 1343
$ws.Range("I89").Value = 827
$ws.Range("J89").Value = 1755.8
$ws.Range("K89").Value = 4135
$ws.Range("L89").Value = 8779
$ws.Range("M89").Value = 1481
$ws.Range("N89").Value = -20011
$ws.Range("H112").Value = 4379.8
$ws.Range("J112").Value = 4379.8
$ws.Range("L112").Value = 13139.4
$ws.Range("N112").Value = -15355.4
$ws.Range("H113").Value = 3091.7778
$ws.Range("I113").Value = 2270
$ws.Range("J113").Value = 3502.6667
$ws.Range("K113").Value = 2270
$ws.Range("L113").Value = 3502.6667
$ws.Range("M113").Value = 984
$ws.Range("N113").Value = -10010.6667
$ws.Range("H121").Value = 770.65717
$ws.Range("I121").Value = 897.5
$ws.Range("J121").Value = 754.29034
$ws.Range("K121").Value = 2692.5
$ws.Range("L121").Value = 2262.87102
$ws.Range("M121").Value = -945.5
$ws.Range("N121").Value = -5756.87102
$ws.Range("H138").Value = 3756.449
$ws.Range("I138").Value = 2178.1667
$ws.Range("J138").Value = 5271.6
$ws.Range("K138").Value = 6534.500100000001
$ws.Range("L138").Value = 15814.8
$ws.Range("M138").Value = -1394.500100000001
$ws.Range("N138").Value = -26094.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12696.181
$ws.Range("I32").Value = 14821.667
$ws.Range("K32").Value = 14821.667
$ws.Range("M32").Value = -14534.667
$ws.Range("H74").Value = 1539.12
$ws.Range("I74").Value = 1440.5483
$ws.Range("J74").Value = 1699.9474
$ws.Range("K74").Value = 1440.5483
$ws.Range("L74").Value = 1699.9474
$ws.Range("M74").Value = -566.5482999999999
$ws.Range("N74").Value = -3447.9474
$ws.Range("H77").Value = 1539.12
$ws.Range("I77").Value = 1440.5483
$ws.Range("J77").Value = 1699.9474
$ws.Range("K77").Value = 7202.7415
$ws.Range("L77").Value = 8499.737000000001
$ws.Range("M77").Value = -2834.7415
$ws.Range("N77").Value = -17235.737

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3397
$ws.Range("I134").Value = 3485.1667
$ws.Range("K134").Value = 10455.5001
$ws.Range("M134").Value = -7920.500100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49760
$ws.Range("H27").Value = 50000
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49808
$ws.Range("H57").Value = 186666.67
$ws.Range("J57").Value = 186666.67
$ws.Range("L57").Value = 186666.67
$ws.Range("N57").Value = -187786.67
$ws.Range("H99").Value = 1964.1428
$ws.Range("I99").Value = 2091.6667
$ws.Range("J99").Value = 1199
$ws.Range("K99").Value = 2091.6667
$ws.Range("L99").Value = 1199
$ws.Range("M99").Value = -593.6667000000002
$ws.Range("N99").Value = -4195
$ws.Range("H126").Value = 1964.1428
$ws.Range("I126").Value = 2091.6667
$ws.Range("J126").Value = 1199
$ws.Range("K126").Value = 6275.000100000001
$ws.Range("L126").Value = 3597
$ws.Range("M126").Value = -3805.000100000001
$ws.Range("N126").Value = -8537
$ws.Range("H134").Value = 1688.6207
$ws.Range("I134").Value = 1502.3334
$ws.Range("J134").Value = 2582.8
$ws.Range("K134").Value = 4507.0002
$ws.Range("L134").Value = 7748.400000000001
$ws.Range("M134").Value = -1972.0002
$ws.Range("N134").Value = -12818.4
$ws.Range("H141").Value = 49600
$ws.Range("J141").Value = 49600
$ws.Range("L141").Value = 49600
$ws.Range("N141").Value = -59960

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 738.15
$ws.Range("I107").Value = 695.075
$ws.Range("K107").Value = 2085.225
$ws.Range("M107").Value = -165.2250000000004
$ws.Range("H131").Value = 2013.1978
$ws.Range("J131").Value = 2389.4385
$ws.Range("L131").Value = 7168.315500000001
$ws.Range("N131").Value = -17248.3155
$ws.Range("H138").Value = 2845.1177
$ws.Range("I138").Value = 915.1667
$ws.Range("J138").Value = 3897.818
$ws.Range("K138").Value = 2745.5001
$ws.Range("L138").Value = 11693.454
$ws.Range("M138").Value = 2394.4999
$ws.Range("N138").Value = -21973.454

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H134").Value = 97757.71000000001
$ws.Range("J134").Value = 97757.71000000001
$ws.Range("L134").Value = 293273.13
$ws.Range("N134").Value = -298343.13

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H68").Value = 2470
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2470
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4898.6875
$ws.Range("I62").Value = 4828.5713
$ws.Range("J62").Value = 4953.222
$ws.Range("K62").Value = 4828.5713
$ws.Range("L62").Value = 4953.222
$ws.Range("M62").Value = -4204.5713
$ws.Range("N62").Value = -6201.222
$ws.Range("H65").Value = 4898.6875
$ws.Range("I65").Value = 4828.5713
$ws.Range("J65").Value = 4953.222
$ws.Range("K65").Value = 24142.8565
$ws.Range("L65").Value = 24766.11
$ws.Range("M65").Value = -21022.8565
$ws.Range("N65").Value = -31006.11
$ws.Range("H69").Value = 8000
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 8000
$ws.Range("N69").Value = -9498
$ws.Range("H72").Value = 8000
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 24000
$ws.Range("N72").Value = -31488
